{"js": "// Map of old multiplication expressions -> new ones, per the commit diff.\nconst replacements = [\n  [\"14\u00d720=\", \"31\u00d738=\"],\n  [\"28\u00d751=\", \"18\u00d734=\"],\n  [\"15\u00d799=\", \"83\u00d790=\"],\n  [\"47\u00d748=\", \"22\u00d714=\"],\n  [\"91\u00d777=\", \"74\u00d763=\"],\n  [\"89\u00d734=\", \"37\u00d735=\"],\n  [\"38\u00d775=\", \"73\u00d746=\"],\n  [\"92\u00d734=\", \"86\u00d714=\"],\n  [\"30\u00d762=\", \"40\u00d776=\"],\n  [\"20\u00d736=\", \"36\u00d772=\"],\n  [\"39\u00d757=\", \"73\u00d788=\"],\n  [\"20\u00d764=\", \"84\u00d754=\"],\n  [\"56\u00d778=\", \"25\u00d721=\"],\n  [\"57\u00d716=\", \"36\u00d711=\"],\n  [\"68\u00d785=\", \"42\u00d736=\"],\n  [\"38\u00d753=\", \"48\u00d770=\"],\n  [\"93\u00d724=\", \"89\u00d796=\"],\n  [\"93\u00d789=\", \"84\u00d777=\"],\n  [\"25\u00d761=\", \"21\u00d770=\"],\n  [\"57\u00d790=\", \"57\u00d788=\"],\n  [\"86\u00d745=\", \"60\u00d742=\"],\n  [\"48\u00d758=\", \"53\u00d779=\"],\n  [\"87\u00d728=\", \"91\u00d764=\"],\n  [\"20\u00d746=\", \"17\u00d747=\"],\n  [\"14\u00d712=\", \"64\u00d767=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each multiplication expression in the table with its updated value,\n# per the commit diff. Every old value is unique in the document, so a simple\n# Find/Replace (ReplaceAll) against the whole document body is safe for each.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"14\u00d720=\"; New = \"31\u00d738=\" },\n    @{ Old = \"28\u00d751=\"; New = \"18\u00d734=\" },\n    @{ Old = \"15\u00d799=\"; New = \"83\u00d790=\" },\n    @{ Old = \"47\u00d748=\"; New = \"22\u00d714=\" },\n    @{ Old = \"91\u00d777=\"; New = \"74\u00d763=\" },\n    @{ Old = \"89\u00d734=\"; New = \"37\u00d735=\" },\n    @{ Old = \"38\u00d775=\"; New = \"73\u00d746=\" },\n    @{ Old = \"92\u00d734=\"; New = \"86\u00d714=\" },\n    @{ Old = \"30\u00d762=\"; New = \"40\u00d776=\" },\n    @{ Old = \"20\u00d736=\"; New = \"36\u00d772=\" },\n    @{ Old = \"39\u00d757=\"; New = \"73\u00d788=\" },\n    @{ Old = \"20\u00d764=\"; New = \"84\u00d754=\" },\n    @{ Old = \"56\u00d778=\"; New = \"25\u00d721=\" },\n    @{ Old = \"57\u00d716=\"; New = \"36\u00d711=\" },\n    @{ Old = \"68\u00d785=\"; New = \"42\u00d736=\" },\n    @{ Old = \"38\u00d753=\"; New = \"48\u00d770=\" },\n    @{ Old = \"93\u00d724=\"; New = \"89\u00d796=\" },\n    @{ Old = \"93\u00d789=\"; New = \"84\u00d777=\" },\n    @{ Old = \"25\u00d761=\"; New = \"21\u00d770=\" },\n    @{ Old = \"57\u00d790=\"; New = \"57\u00d788=\" },\n    @{ Old = \"86\u00d745=\"; New = \"60\u00d742=\" },\n    @{ Old = \"48\u00d758=\"; New = \"53\u00d779=\" },\n    @{ Old = \"87\u00d728=\"; New = \"91\u00d764=\" },\n    @{ Old = \"20\u00d746=\"; New = \"17\u00d747=\" },\n    @{ Old = \"14\u00d712=\"; New = \"64\u00d767=\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
